# Se arregla el flujo de archivos PREI y SAI
# Update column O ("Factura") for rows that previously held "sin match"
# placeholders, assigning them their matched invoice/folio numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("df_altas")

$updates = @{
    4  = "P-6570"
    10 = "P-6572"
    11 = "P-6573"
    12 = "P-6574"
    14 = "P-6575"
    22 = "P-6576"
    24 = "P-6577"
    25 = "P-6578"
    30 = "P-6579"
    33 = "P-6580"
    36 = "P-6581"
    40 = "P-6582"
    41 = "P-6583"
    46 = "P-6584"
}

foreach ($row in $updates.Keys) {
    $ws.Range("O$row").Value = $updates[$row]
}
